$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# New forecast values (Amazon Mean / P70 / P80 / P90) after removing Auto Arima,
# keyed by row number for columns D, E, F, G.
$data = @{
    2  = @(145, 165, 182, 206)
    3  = @(110, 128, 145, 169)
    4  = @(108, 126, 141, 165)
    5  = @(107, 126, 142, 167)
    6  = @(107, 127, 145, 173)
    7  = @(106, 125, 142, 168)
    8  = @(107, 126, 145, 173)
    9  = @(107, 126, 145, 173)
    10 = @(105, 124, 141, 168)
    11 = @(105, 124, 141, 169)
    13 = @(105, 125, 144, 174)
    14 = @(102, 121, 139, 168)
    15 = @(96, 116, 135, 166)
    16 = @(94, 113, 132, 161)
    17 = @(93, 112, 131, 160)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
    $ws.Range("F$row").Value = $vals[2]
    $ws.Range("G$row").Value = $vals[3]
}
